$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values (B:G) for rows 2-11, regenerated sval data filtering save games.
$data = @{
    2  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248)
    3  = @(1.445647641019636, 0.04103571897497393, 0.7210945179870265, 0.5333859586016987, 0, 2.741163836583335)
    4  = @(0.1169995834814548, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 0, 1.104883657715537)
    5  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248)
    6  = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 13.86384647080068, 1, 21.98653043760045)
    7  = @(0.01253208636536152, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 0, 0.7365606309410384)
    8  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
    9  = @(0.1169995834814548, 0.3048912486333797, 0.7210945179870265, 13.86384647080068, 1, 15.00683182090255)
    10 = @(0.6545652718822623, 1.626987699542094, 189.6080260415259, 2797.565817734744, 0, 2989.455396747694)
    11 = @(3.272327238179451, 9.983522426115931, 18.71679738969934, 13.86384647080068, 1, 45.8364935247954)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
